$wb = $excel.ActiveWorkbook

# --- BattleStage sheet: add a new "EntityList" field (column G) -----------
$ws = $wb.Worksheets.Item("BattleStage")

$ws.Range("G2").Value = "All"
$ws.Range("G3").Value = "EntityList"
$ws.Range("G4").Value = "List<StageActor{Uid:Int32,Rid:Int32,Level:Int32,PosX:Int32,PosY:Int32}>"
$ws.Range("G8").Value = "101,1,4,5,3 | 102,1,4,5,1"

# NB: Excel quantises ColumnWidth to whole pixels, so the closest
# achievable stored width to the source file's 61.375 is reached by
# feeding it a char-width of 60.7 (-> ~61.43 once rendered).
$ws.Columns.Item(7).ColumnWidth = 60.7

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Update selections on the other sheets ---------------------------------
$wsScene = $wb.Worksheets.Item("BattleScene")
[void]$wsScene.Range("D8").Select()

$wsTerrain = $wb.Worksheets.Item("Terrain")
[void]$wsTerrain.Range("E8").Select()

# --- BattleStage becomes the active tab (selected last) --------------------
$ws.Activate()
[void]$ws.Range("G8").Select()
